# Applies the Fri Mar 31 14:47:13 UTC 2023 "cryptos list" refresh:
#   - updates Price (D) / Volume(1h) (E) columns for each coin row
#   - row 26/28 additionally swap Coin name + Link (ranking order changed)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values look numeric (e.g. "0.9988"); Excel would silently coerce a plain
# Value assignment into a float and drop the original text formatting. Force the
# cell to Text first, write the literal string, then clear the temporary format so
# the cell is left exactly as it started (no left-over "@" number format).
foreach ($addr in @(
        "D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15",
        "D16", "D18", "D19", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29",
        "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41",
        "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51"
    )) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '28.525.99'
$ws.Range("D3").Value = '1.840.99'
$ws.Range("D4").Value = '0.9988'
$ws.Range("D5").Value = '319.10'
$ws.Range("D6").Value = '0.9991'
$ws.Range("D7").Value = '0.5337'
$ws.Range("D8").Value = '0.4016'
$ws.Range("D9").Value = '0.07597'
$ws.Range("D10").Value = '41.89'
$ws.Range("D11").Value = '1.112'
$ws.Range("D12").Value = '6.329'
$ws.Range("D13").Value = '7.626'
$ws.Range("D14").Value = '0.9992'
$ws.Range("D15").Value = '20.84'
$ws.Range("D16").Value = '1.833.98'
$ws.Range("D18").Value = '0.00001074'
$ws.Range("D19").Value = '0.06601'
$ws.Range("D22").Value = '6.070'
$ws.Range("D23").Value = '28.538.18'
$ws.Range("D24").Value = '11.23'
$ws.Range("D25").Value = '2.104'
$ws.Range("D26").Value = '2.464'
$ws.Range("D27").Value = '157.10'
$ws.Range("D28").Value = '20.69'
$ws.Range("D29").Value = '2.044.10'
$ws.Range("D30").Value = '123.98'
$ws.Range("D31").Value = '1.124'
$ws.Range("D33").Value = '5.709'
$ws.Range("D34").Value = '3.664'
$ws.Range("D35").Value = '0.07229'
$ws.Range("D36").Value = '0.2261'
$ws.Range("D37").Value = '5.271'
$ws.Range("D38").Value = '0.02349'
$ws.Range("D39").Value = '8.829'
$ws.Range("D40").Value = '11.37'
$ws.Range("D41").Value = '0.6295'
$ws.Range("D42").Value = '1.206'
$ws.Range("D43").Value = '1.413'
$ws.Range("D44").Value = '0.9986'
$ws.Range("D45").Value = '13.53'
$ws.Range("D46").Value = '3.712'
$ws.Range("D47").Value = '0.5853'
$ws.Range("D48").Value = '126.04'
$ws.Range("D49").Value = '1.982'
$ws.Range("D51").Value = '0.06928'

foreach ($addr in @(
        "D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15",
        "D16", "D18", "D19", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29",
        "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41",
        "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51"
    )) {
    $ws.Range($addr).ClearFormats()
}

# Coin name / link swaps and Volume(1h) text updates (plain text, no coercion risk)
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("E3").Value = '  +2.40%  '
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("E6").Value = '  -0.26%  '
$ws.Range("E7").Value = '  -1.61%  '
$ws.Range("E8").Value = '  +6.13%  '
$ws.Range("E9").Value = '  +1.25%  '
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("E12").Value = '  +2.75%  '
$ws.Range("E13").Value = '  +4.45%  '
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("E16").Value = '  +2.25%  '
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("E18").Value = '  +0.68%  '
$ws.Range("E19").Value = '  +1.51%  '
$ws.Range("E20").Value = '  +1.63%  '
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("E22").Value = '  +1.92%  '
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("E24").Value = '  +1.11%  '
$ws.Range("E25").Value = '  +1.48%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("E26").Value = '  +6.06%  '
$ws.Range("E27").Value = '  -1.56%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("E28").Value = '  +1.16%  '
$ws.Range("E29").Value = '  +2.16%  '
$ws.Range("E30").Value = '  +0.98%  '
$ws.Range("E31").Value = '  +1.21%  '
$ws.Range("E32").Value = '  +4.29%  '
$ws.Range("E33").Value = '  +1.69%  '
$ws.Range("E34").Value = '  +0.33%  '
$ws.Range("E35").Value = '  +11.85%  '
$ws.Range("E36").Value = '  -0.76%  '
$ws.Range("E37").Value = '  +4.72%  '
$ws.Range("E38").Value = '  +2.09%  '
$ws.Range("E39").Value = '  +2.42%  '
$ws.Range("E40").Value = '  +1.60%  '
$ws.Range("E41").Value = '  +1.43%  '
$ws.Range("E42").Value = '  +1.02%  '
$ws.Range("E43").Value = '  -2.57%  '
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("E45").Value = '  +1.68%  '
$ws.Range("E46").Value = '  +0.70%  '
$ws.Range("E47").Value = '  +0.41%  '
$ws.Range("E48").Value = '  -0.94%  '
$ws.Range("E49").Value = '  +1.62%  '
$ws.Range("E50").Value = '  -0.76%  '
$ws.Range("E51").Value = '  +0.71%  '
